$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells with the same style as existing headers (row 1, column A-AC use style index 1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the style from an existing header cell (e.g. AC1) to the new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in Wins/Losses/Ties values for every data row (2 through 54)
for ($r = 2; $r -le 54; $r++) {
    $ws.Cells.Item($r, 30).Value = 51   # AD = column 30
    $ws.Cells.Item($r, 31).Value = 111  # AE = column 31
    $ws.Cells.Item($r, 32).Value = 0    # AF = column 32
}
